$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "clear option" on the sign-up sheet let anyone wipe rows without a
# password, and a couple of joke/fake course entries ("我只是来求课的TvT",
# "我是雷锋我就是来出课的") had been submitted in rows 31/32 (sheet rows
# 47-48). Remove those fake rows, keeping their formatting intact.
$ws.Range("A47:P48").ClearContents()

# Leave the selection where the edit was made.
$ws.Range("R40").Select()
